$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.8629528957341
$ws.Range("D2").Value = 8.942448215102498
$ws.Range("E2").Value = 14.29024043334591
$ws.Range("F2").Value = 35.49363219902544
$ws.Range("G2").Value = 3.674716527338106
$ws.Range("I2").Value = 28.43855552673963
$ws.Range("J2").Value = 10.78638394133868
$ws.Range("K2").Value = 10.03716883476748
$ws.Range("L2").Value = 10.14703833950733
$ws.Range("O2").Value = 26.91738179466552

$ws.Range("B3").Value = 13.71410000655744
$ws.Range("D3").Value = 8.944783443083029
$ws.Range("E3").Value = 14.32457908004933
$ws.Range("F3").Value = 35.59502380556452
$ws.Range("G3").Value = 3.676589619671166
$ws.Range("I3").Value = 28.56434485224948
$ws.Range("J3").Value = 10.80684356715101
$ws.Range("K3").Value = 9.762976812841375
$ws.Range("L3").Value = 10.10566498267973
$ws.Range("O3").Value = 27.00805034755464

$ws.Range("B4").Value = 13.62381312636515
$ws.Range("D4").Value = 8.947202664733933
$ws.Range("E4").Value = 14.34705506051124
$ws.Range("F4").Value = 35.66475591643341
$ws.Range("G4").Value = 3.677801436754113
$ws.Range("I4").Value = 28.64608231189364
$ws.Range("J4").Value = 10.82008088483117
$ws.Range("K4").Value = 9.591599268351661
$ws.Range("L4").Value = 10.08144789431045
$ws.Range("O4").Value = 27.06903210795875

$ws.Range("B5").Value = 13.58733578573416
$ws.Range("D5").Value = 8.948437097012734
$ws.Range("E5").Value = 14.35656493969131
$ws.Range("F5").Value = 35.69504952568217
$ws.Range("G5").Value = 3.678310832735795
$ws.Range("I5").Value = 28.68052474115756
$ws.Range("J5").Value = 10.82564542460316
$ws.Range("K5").Value = 9.521100875474929
$ws.Range("L5").Value = 10.07188413136261
$ws.Range("O5").Value = 27.09521638203076

$ws.Range("B6").Value = 13.58129884562683
$ws.Range("D6").Value = 8.948657109509115
$ws.Range("E6").Value = 14.35816525280937
$ws.Range("F6").Value = 35.7001930430493
$ws.Range("G6").Value = 3.678396359392875
$ws.Range("I6").Value = 28.68631240395978
$ws.Range("J6").Value = 10.82657970765976
$ws.Range("K6").Value = 9.509357793143588
$ws.Range("L6").Value = 10.07031467712126
$ws.Range("O6").Value = 27.09964476125155

$ws.Range("B7").Value = 13.62331985596005
$ws.Range("D7").Value = 8.947218305201934
$ws.Range("E7").Value = 14.34718189297951
$ws.Range("F7").Value = 35.66515687004336
$ws.Range("G7").Value = 3.677808243537588
$ws.Range("I7").Value = 28.64654222187172
$ws.Range("J7").Value = 10.82015524027694
$ws.Range("K7").Value = 9.59065104151934
$ws.Range("L7").Value = 10.0813176710755
$ws.Range("O7").Value = 27.06937984046473

$ws.Range("B8").Value = 13.81142010146824
$ws.Range("D8").Value = 8.943049377041902
$ws.Range("E8").Value = 14.30179194482514
$ws.Range("F8").Value = 35.52703837477804
$ws.Range("G8").Value = 3.675349586261448
$ws.Range("I8").Value = 28.48099420037105
$ws.Range("J8").Value = 10.79329861139712
$ws.Range("K8").Value = 9.943312463051692
$ws.Range("L8").Value = 10.13253062847771
$ws.Range("O8").Value = 26.94754099510118

$ws.Range("B9").Value = 14.18747270281207
$ws.Range("D9").Value = 8.942656480635094
$ws.Range("E9").Value = 14.22379383464124
$ws.Range("F9").Value = 35.31564413072081
$ws.Range("G9").Value = 3.67101578466635
$ws.Range("I9").Value = 28.19200511213835
$ws.Range("J9").Value = 10.74596668558367
$ws.Range("K9").Value = 10.60685246852416
$ws.Range("L9").Value = 10.24207000559595
$ws.Range("O9").Value = 26.75083174205186

$ws.Range("B10").Value = 14.46596879952408
$ws.Range("D10").Value = 8.947061782047582
$ws.Range("E10").Value = 14.17315622092067
$ws.Range("F10").Value = 35.19674087252421
$ws.Range("G10").Value = 3.668125942569788
$ws.Range("I10").Value = 28.00131468088281
$ws.Range("J10").Value = 10.71441246028514
$ws.Range("K10").Value = 11.07231816356467
$ws.Range("L10").Value = 10.32768109231356
$ws.Range("O10").Value = 26.63214772137665

$ws.Range("B11").Value = 14.59268458502804
$ws.Range("D11").Value = 8.950073359612887
$ws.Range("E11").Value = 14.15155810025895
$ws.Range("F11").Value = 35.15058419654195
$ws.Range("G11").Value = 3.666874510594266
$ws.Range("I11").Value = 27.9192400465124
$ws.Range("J11").Value = 10.70075032327397
$ws.Range("K11").Value = 11.27837328422968
$ws.Range("L11").Value = 10.36764707517952
$ws.Range("O11").Value = 26.58378647441403

$ws.Range("B12").Value = 14.64063442439151
$ws.Range("D12").Value = 8.951357550864952
$ws.Range("E12").Value = 14.14358539256637
$ws.Range("F12").Value = 35.13424862282314
$ws.Range("G12").Value = 3.666409660593614
$ws.Range("I12").Value = 27.88883080327494
$ws.Range("J12").Value = 10.69567584448621
$ws.Range("K12").Value = 11.35551605703915
$ws.Range("L12").Value = 10.38291956023395
$ws.Range("O12").Value = 26.56628435649182

$ws.Range("B13").Value = 14.63030971641653
$ws.Range("D13").Value = 8.951074601388738
$ws.Range("E13").Value = 14.14529330452456
$ws.Range("F13").Value = 35.13771591570468
$ws.Range("G13").Value = 3.666509372991796
$ws.Range("I13").Value = 27.89535016805547
$ws.Range("J13").Value = 10.69676432479405
$ws.Range("K13").Value = 11.33894239208632
$ws.Range("L13").Value = 10.37962434817795
$ws.Range("O13").Value = 26.57001764246348

$ws.Range("B14").Value = 14.59663036021647
$ws.Range("D14").Value = 8.950176135991001
$ws.Range("E14").Value = 14.15089805522068
$ws.Range("F14").Value = 35.14921734094607
$ws.Range("G14").Value = 3.666836086177445
$ws.Range("I14").Value = 27.91672482529457
$ws.Range("J14").Value = 10.70033085971595
$ws.Range("K14").Value = 11.28473798636821
$ws.Range("L14").Value = 10.36890083119933
$ws.Range("O14").Value = 26.58233029604178

$ws.Range("B15").Value = 14.57599513511468
$ws.Range("D15").Value = 8.949644489253833
$ws.Range("E15").Value = 14.15435793958633
$ws.Range("F15").Value = 35.15641120304306
$ws.Range("G15").Value = 3.66703738343275
$ws.Range("I15").Value = 27.92990472522066
$ws.Range("J15").Value = 10.70252835501635
$ws.Range("K15").Value = 11.25141895437828
$ws.Range("L15").Value = 10.36235011127636
$ws.Range("O15").Value = 26.58997785184655

$ws.Range("B16").Value = 14.45768485319795
$ws.Range("D16").Value = 8.946885149334015
$ws.Range("E16").Value = 14.17459657443543
$ws.Range("F16").Value = 35.19991714720351
$ws.Range("G16").Value = 3.668208994008697
$ws.Range("I16").Value = 28.00677235645162
$ws.Range("J16").Value = 10.71531920109398
$ws.Range("K16").Value = 11.0587315710272
$ws.Range("L16").Value = 10.32508903779391
$ws.Range("O16").Value = 26.63542166544155

$ws.Range("B17").Value = 14.38508453050104
$ws.Range("D17").Value = 8.945449673219951
$ws.Range("E17").Value = 14.18737996610856
$ws.Range("F17").Value = 35.22864012893637
$ws.Range("G17").Value = 3.668943887594232
$ws.Range("I17").Value = 28.05512381651168
$ws.Range("J17").Value = 10.72334290538214
$ws.Range("K17").Value = 10.93901696276289
$ws.Range("L17").Value = 10.30248604629056
$ws.Range("O17").Value = 26.6647429532391

$ws.Range("B18").Value = 14.34333150794043
$ws.Range("D18").Value = 8.944718934342845
$ws.Range("E18").Value = 14.19486794696303
$ws.Range("F18").Value = 35.24590740785779
$ws.Range("G18").Value = 3.669372527565439
$ws.Range("I18").Value = 28.08337405451839
$ws.Range("J18").Value = 10.72802309175045
$ws.Range("K18").Value = 10.86962859021887
$ws.Range("L18").Value = 10.28958203908642
$ws.Range("O18").Value = 26.6821374786535

$ws.Range("B19").Value = 14.3291966533031
$ws.Range("D19").Value = 8.944487852531104
$ws.Range("E19").Value = 14.19742650798636
$ws.Range("F19").Value = 35.25188196569427
$ws.Range("G19").Value = 3.669518680719273
$ws.Range("I19").Value = 28.09301466277984
$ws.Range("J19").Value = 10.72961892781602
$ws.Range("K19").Value = 10.84604573733128
$ws.Range("L19").Value = 10.28522981980593
$ws.Range("O19").Value = 26.68811788462914

$ws.Range("B20").Value = 14.39281274677012
$ws.Range("D20").Value = 8.945592667296788
$ws.Range("E20").Value = 14.18600515256376
$ws.Range("F20").Value = 35.22550523342399
$ws.Range("G20").Value = 3.668865041576683
$ws.Range("I20").Value = 28.04993121596887
$ws.Range("J20").Value = 10.72248202733471
$ws.Range("K20").Value = 10.95181632854947
$ws.Range("L20").Value = 10.30488223351594
$ws.Range("O20").Value = 26.66156681695077

$ws.Range("B21").Value = 14.60652404805404
$ws.Range("D21").Value = 8.950436144306174
$ws.Range("E21").Value = 14.14924621745489
$ws.Range("F21").Value = 35.1458080579575
$ws.Range("G21").Value = 3.666739877596523
$ws.Range("I21").Value = 27.91042837791157
$ws.Range("J21").Value = 10.69928059653472
$ws.Range("K21").Value = 11.30068368519893
$ws.Range("L21").Value = 10.37204690587552
$ws.Range("O21").Value = 26.57869174029213

$ws.Range("B22").Value = 14.74597726785426
$ws.Range("D22").Value = 8.954439092659872
$ws.Range("E22").Value = 14.12642274988436
$ws.Range("F22").Value = 35.10038418637834
$ws.Range("G22").Value = 3.665403632771919
$ws.Range("I22").Value = 27.82316373254231
$ws.Range("J22").Value = 10.68469441040266
$ws.Range("K22").Value = 11.52349832928536
$ws.Range("L22").Value = 10.41674452250343
$ws.Range("O22").Value = 26.52925738493009

$ws.Range("B23").Value = 14.67158107300226
$ws.Range("D23").Value = 8.952226400179898
$ws.Range("E23").Value = 14.13849441490401
$ws.Range("F23").Value = 35.12401746214522
$ws.Range("G23").Value = 3.666112006631897
$ws.Range("I23").Value = 27.86938123742278
$ws.Range("J23").Value = 10.6924266541478
$ws.Range("K23").Value = 11.40507368246812
$ws.Range("L23").Value = 10.39281810277256
$ws.Range("O23").Value = 26.55520811583854

$ws.Range("B24").Value = 14.38931886126644
$ws.Range("D24").Value = 8.945527725138858
$ws.Range("E24").Value = 14.18662627388588
$ws.Range("F24").Value = 35.22692017102044
$ws.Range("G24").Value = 3.6689006687328
$ws.Range("I24").Value = 28.05227738157081
$ws.Range("J24").Value = 10.7228710207855
$ws.Range("K24").Value = 10.94603148319802
$ws.Range("L24").Value = 10.30379863393286
$ws.Range("O24").Value = 26.66300107455356

$ws.Range("B25").Value = 14.08519528196635
$ws.Range("D25").Value = 8.941934597906721
$ws.Range("E25").Value = 14.24372024117176
$ws.Range("F25").Value = 35.36644897712846
$ws.Range("G25").Value = 3.672136307951125
$ws.Range("I25").Value = 28.26637859107042
$ws.Range("J25").Value = 10.75820343182761
$ws.Range("K25").Value = 10.43088526653685
$ws.Range("L25").Value = 10.21150375916381
$ws.Range("O25").Value = 26.79951645421159
